$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D" = -0.241
    "G" = -0.007407407407407407
    "H" = -0.01703703703703703
    "I" = -0.2162962962962963
    "J" = -0.2162962962962963
    "K" = -0.272
    "L" = -0.2014814814814815
    "U" = 0.002
    "V" = 0.001058201058201058
    "W" = -0.09927007299270073
    "X" = 0.08399743193505964
    "Y" = -0.1832675049277604
    "Z" = 0.3613490364025695
    "AA" = -0.07815845824411133
    "AB" = 0.07767800168920949
    "AC" = -0.1558364599333208
    "AD" = 0.323
    "AF" = 0.323
    "AG" = 0.321
    "AH" = 0.1459557162223226
    "AI" = 0.1306105944197331
    "AJ" = 0.1451831750339213
    "AK" = 0.1299069202751922
    "AL" = 0.1
    "AM" = 0.1
    "AN" = -1.386266094420601
    "AO" = -2.919999999999999
    "AP" = -1.377682403433476
    "AQ" = -2.919999999999999
}

foreach ($row in 2..3) {
    foreach ($col in $updates.Keys) {
        $ws.Range("$col$row").Value = $updates[$col]
    }
}
